$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Activation date: 01/01/2022 -> 01/01/2025
# ------------------------------------------------------------------
$d.Content.Find.Execute(
    "Ativação: 01/01/2022", $true, $false, $false, $false, $false,
    $true, 1, $false, "Ativação: 01/01/2025", 2)

# ------------------------------------------------------------------
# 2. English "Objetivos" paragraph gains a trailing period.
# ------------------------------------------------------------------
$d.Content.Find.Execute(
    "Understand the theoretical bases of qualitative analytical chemistry of environmental interest",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Understand the theoretical bases of qualitative analytical chemistry of environmental interest.",
    2)

# ------------------------------------------------------------------
# 3. Insert four new teachers before the existing one in the
#    "Docente(s) Responsável(eis)" bullet list.
# ------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("7455355 - Robson da Silva Rocha", $true, $false, $false, $false, $false,
                 $true, 1, $false, "", 0)
$r.Collapse(1)
$nl = [char]11
$r.InsertBefore("7043088 - Ana Karine Furtado de Carvalho" + $nl)
$r.Collapse(0)
$r.InsertBefore("7926291 - Célia Regina Tomachuk dos Santos Catuogno" + $nl)
$r.Collapse(0)
$r.InsertBefore("4893449 - Débora Souza Alvim" + $nl)
$r.Collapse(0)
$r.InsertBefore("8855158 - Morun Bernardino Neto" + $nl)

# ------------------------------------------------------------------
# 4. Trim the Portuguese "Programa resumido" paragraph.
# ------------------------------------------------------------------
$d.Content.Find.Execute(
    " Análise de sólidos, partículas, sedimentos. Estudos de amostras de importância ambiental.",
    $true, $false, $false, $false, $false,
    $true, 1, $false, "", 2)

# ------------------------------------------------------------------
# 5. Trim the English "Programa resumido" paragraph.
# ------------------------------------------------------------------
$d.Content.Find.Execute(
    " Analysis of solids, particles, sediments. Studies of samples of environmental importance.",
    $true, $false, $false, $false, $false,
    $true, 1, $false, "", 2)

# ------------------------------------------------------------------
# 6. Trim the Portuguese "Programa" paragraph (two removals).
# ------------------------------------------------------------------
$d.Content.Find.Execute(
    "- Análise de sólidos, partículas, sedimentos.",
    $true, $false, $false, $false, $false,
    $true, 1, $false, "", 2)

$d.Content.Find.Execute(
    "- Análise gravimétrica: fundamentos e formação de precipitados.- Análises dos principais cátions e ânions em amostras conhecidas e desconhecidas para os alunos- Análise de metais em solo, água ou outras amostras ambientais importantes",
    $true, $false, $false, $false, $false,
    $true, 1, $false, "", 2)

# ------------------------------------------------------------------
# 7. Trim the English "Programa" paragraph (two removals).
# ------------------------------------------------------------------
$d.Content.Find.Execute(
    "- Analysis of solids, particles, sediments.",
    $true, $false, $false, $false, $false,
    $true, 1, $false, "", 2)

$d.Content.Find.Execute(
    "- Gravimetric analysis: fundamentals and precipitate formation.- Analysis of the main cations and anions in known and unknown samples for students- Analysis of metals in soil, water or other important environmental samples",
    $true, $false, $false, $false, $false,
    $true, 1, $false, "", 2)
